$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Range("B8")
$cell.Formula = "'= addAll1(null, null); " + [char]34 + "Hello" + [char]34 + ";"
[void]$ws.Range("D10").Select()
